# Automatically add the path of the addin file to the Python path so that
# .py files are found easily.
#
# The three "bootstrap" example rows that used to show the user typing
# `import ExPy;`, `import sys;` and `sys.path.append("F:\\")` by hand are
# no longer needed (the add-in now does this for you), so remove those
# demo cells from both example sheets. Everything else (the D/E "eval"
# examples on RegEx, and the numpy examples) is left in place.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- "RegEx" sheet: remove the ExPyScript(...) demo rows that referenced
# the bootstrap strings (A2:B2, A3:B3, A4:B4). ---
$ws1.Range("A2:B2").ClearContents()
$ws1.Range("A3:B3").ClearContents()
$ws1.Range("A4:B4").ClearContents()

# --- "use numpy" sheet: remove the same kind of demo rows (the whole of
# row 1, plus B2:C2 and B3:C3). ---
$ws2.Range("B1:C1").ClearContents()
$ws2.Range("B2:C2").ClearContents()
$ws2.Range("B3:C3").ClearContents()

# Reset the view state: "use numpy" is no longer the active tab and loses
# its custom selection, while "RegEx" becomes active with D5 selected.
$ws2.Range("A1").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("D5").Select() | Out-Null
